# AvisVirement.docx update ("bug fixed - inhance UX")
#   - DAS reference year   "2020"      -> "2021"   (appears twice in the doc)
#   - Amount                "3 767 400,00" -> "3 993 300,00"
#   - Payment month/year    "Decembre 2020" -> "Janvier 2021"
#
# Most of the text lives in the normal document story and is reachable with
# a plain Find/Replace. One of the two "2020" occurrences, however, sits
# inside a legacy VML drawing (a text box drawn with <w:pict>/<v:textbox>)
# that Word's navigable collections (Content, Paragraphs, Characters,
# Shapes, StoryRanges, ...) do not expose in this runtime. To reach that
# text too we round-trip the whole package through WordOpenXML/InsertXML,
# doing plain text substitutions on the serialized XML before handing it
# back - this is the one reliable way to touch text that lives inside the
# drawing while leaving every other part of the package untouched.

$d = $word.ActiveDocument

# --- Edits reachable through the normal Word object model -----------------
$d.Content.Find.Execute("3 767 400,00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3 993 300,00", 2) | Out-Null

$d.Content.Find.Execute("Decembre", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Janvier", 2) | Out-Null

$d.Content.Find.Execute("2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2021", 2) | Out-Null

# --- Remaining edit: the "2020" that lives inside the VML text box --------
$xml = $d.Content.WordOpenXML
$patched = $xml.Replace("<w:t>2020</w:t>", "<w:t>2021</w:t>")
if ($patched -ne $xml) {
    $d.Content.InsertXML($patched)
}
